$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")
$ws.Activate()

# Finalize row 21: previously-blank trailing cells (B:K, M) become the
# literal string "nan" (matches the pattern used by every other row).
$ws.Range("B21").Value = "nan"
$ws.Range("C21").Value = "nan"
$ws.Range("D21").Value = "nan"
$ws.Range("E21").Value = "nan"
$ws.Range("F21").Value = "nan"
$ws.Range("G21").Value = "nan"
$ws.Range("H21").Value = "nan"
$ws.Range("I21").Value = "nan"
$ws.Range("J21").Value = "nan"
$ws.Range("K21").Value = "nan"
$ws.Range("M21").Value = "nan"

# Append the new event row (row 22) for Card20.
# Force A22/L22 to be stored as text (not auto-converted to numbers/dates)
# the same way the other "card"/"Date" columns in this sheet are, then
# restore the default "Normal" style so no stray number-format style gets
# attached to the cell.
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "20"
$ws.Range("A22").Style = "Normal"

$ws.Range("L22").NumberFormat = "@"
$ws.Range("L22").Value = "26\2\2025"
$ws.Range("L22").Style = "Normal"

$ws.Range("N22").Value = "تم سن الفلاتس وتغيير الجرائد الخلفيه (1_5_8)"
$ws.Range("O22").Value = "الخبير"
